$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 4 (end of game purpose was missing)
$ws.Range("E4").Value = "End of Game"

# Shorten file name entries for rows 5 and 6
$ws.Range("A5").Value = "Startup"
$ws.Range("A6").Value = "Object Select"

# Add new row 13 for "Collection Chirp" sound
$ws.Range("A13").Value = "Collection Chirp"
$ws.Range("B13").Value = "wav"
$ws.Range("C13").Value = "Jarryd"
$ws.Range("D13").Value = "n/a "
$ws.Range("E13").Value = "used when collecting shards"

# Update selection/view state
$ws.Range("C10").Select()
